# Fix the "Create Teams" sample data rows:
# the planner group ZP1/ZPx confusion meant the sample data was built
# around the ZP1 example; refresh it to a ZPQ (non zp1) example instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Create Teams")

# Update the example data rows (values are written in this specific
# order so that the regenerated shared-string table lines up with the
# canonical file: new strings are interned in first-seen order).
$ws.Range("A2").Value = "0-BR-CEN-01"
$ws.Range("B2").Value = "ETIN0001"
$ws.Range("C2").Value = "ZPQ"
$ws.Range("D2").Value = "B020"
$ws.Range("A3").Value = "0-BR-CEN-02"
$ws.Range("E2").Value = "PREDIRE"
$ws.Range("B3").Value = "ETIN0001"
$ws.Range("C3").Value = "ZPQ"
$ws.Range("D3").Value = "B020"
$ws.Range("E3").Value = "PREDIRE"

# Widen column E (Contractor) to fit the new "PREDIRE" sample text,
# compensating for the engine's internal pixel-rounding so the stored
# OOXML column width lands exactly on 41.
$ws.Columns.Item(5).ColumnWidth = (41 - 5/6)

# Move the saved cursor/selection to E9, matching the author's last
# position when they saved the workbook.
[void]$ws.Range("E9").Select()
